# Updated cryptos list on Fri Apr 28 14:48:32 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'29.203.26"
$ws.Range("E2").Value = "  +0.55%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'1.895.08"
$ws.Range("E3").Value = "  +0.40%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'322.41"
$ws.Range("E5").Value = "  -2.53%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.14%  "

# Row 7 - XRP
$ws.Range("D7").Value = "'0.4695"
$ws.Range("E7").Value = "  +2.31%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "'0.4014"
$ws.Range("E8").Value = "  -2.17%  "

# Row 9 - OKB
$ws.Range("D9").Value = "'47.35"

# Row 10 - Dogecoin
$ws.Range("D10").Value = "'0.07989"
$ws.Range("E10").Value = "  +0.55%  "

# Row 11 - Polygon
$ws.Range("D11").Value = "'0.9927"
$ws.Range("E11").Value = "  -0.39%  "

# Row 12 - Solana
$ws.Range("D12").Value = "'22.37"
$ws.Range("E12").Value = "  +2.47%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "'1.891.09"
$ws.Range("E13").Value = "  -0.63%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "'5.847"

# Row 15 - Chainlink
$ws.Range("D15").Value = "'7.028"
$ws.Range("E15").Value = "  -0.61%  "

# Row 16 / 17 swap: Litecoin <-> BinanceUSD
$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.20%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'88.82"
$ws.Range("E17").Value = "  +0.28%  "

# Row 18 - TRON
$ws.Range("D18").Value = "'0.06621"
$ws.Range("E18").Value = "  +1.07%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "'0.00001025"
$ws.Range("E19").Value = "  +0.08%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "'17.41"
$ws.Range("E20").Value = "  -0.06%  "

# Row 21 - Dai
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.08%  "

# Row 22 - WrappedBTC
$ws.Range("D22").Value = "'29.220.57"
$ws.Range("E22").Value = "  +0.54%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "'5.491"
$ws.Range("E23").Value = "  +1.04%  "

# Row 25 - Toncoin
$ws.Range("E25").Value = "  -0.18%  "

# Row 26 - WrappedliquidstakedEther2.0
$ws.Range("D26").Value = "'2.136.40"
$ws.Range("E26").Value = "  +0.60%  "

# Row 27 - Monero
$ws.Range("D27").Value = "'154.73"
$ws.Range("E27").Value = "  -1.05%  "

# Row 28 - EthereumClassic
$ws.Range("D28").Value = "'19.59"
$ws.Range("E28").Value = "  +0.22%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "'6.070"
$ws.Range("E29").Value = "  +10.28%  "

# Row 30 - LidoDAOToken
$ws.Range("D30").Value = "'2.073"
$ws.Range("E30").Value = "  -0.66%  "

# Row 31 - BitcoinCash
$ws.Range("D31").Value = "'117.06"
$ws.Range("E31").Value = "  -0.17%  "

# Row 32 - ImmutableX
$ws.Range("D32").Value = "'1.047"
$ws.Range("E32").Value = "  +1.22%  "

# Row 33 - Stellar
$ws.Range("E33").Value = "  +1.07%  "

# Row 34 - ARBITRUM
$ws.Range("E34").Value = "  -1.18%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "'3.543"
$ws.Range("E35").Value = "  +0.62%  "

# Row 36 - Filecoin
$ws.Range("D36").Value = "'5.334"
$ws.Range("E36").Value = "  +0.87%  "

# Row 37 - Hedera
$ws.Range("D37").Value = "'0.06062"
$ws.Range("E37").Value = "  +0.17%  "

# Row 39 - TrustWalletToken
$ws.Range("D39").Value = "'1.175"
$ws.Range("E39").Value = "  +0.14%  "

# Row 40 - FraxShare
$ws.Range("E40").Value = "  -3.70%  "

# Row 41 - TheSandbox
$ws.Range("D41").Value = "'0.5800"
$ws.Range("E41").Value = "  +0.32%  "

# Row 42 - Algorand
$ws.Range("D42").Value = "'0.1821"
$ws.Range("E42").Value = "  -0.03%  "

# Row 43 - RenderToken
$ws.Range("D43").Value = "'2.472"
$ws.Range("E43").Value = "  +8.13%  "

# Row 44 - Aptos
$ws.Range("D44").Value = "'10.01"
$ws.Range("E44").Value = "  -0.77%  "

# Row 45 - WEMIXToken
$ws.Range("D45").Value = "'1.272"
$ws.Range("E45").Value = "  +1.11%  "

# Row 46 - Cronos
$ws.Range("D46").Value = "'0.07677"
$ws.Range("E46").Value = "  +2.50%  "

# Row 47 - EnergySwap
$ws.Range("D47").Value = "'12.16"
$ws.Range("E47").Value = "  +1.77%  "

# Row 48 - Decentraland
$ws.Range("D48").Value = "'0.5463"
$ws.Range("E48").Value = "  +0.23%  "

# Row 49 - NEARProtocol
$ws.Range("D49").Value = "'1.896"
$ws.Range("E49").Value = "  -0.37%  "

# Row 50 - Quant
$ws.Range("D50").Value = "'113.50"
$ws.Range("E50").Value = "  +1.55%  "

# Row 51 - Elrond
$ws.Range("D51").Value = "'43.92"
$ws.Range("E51").Value = "  -0.83%  "
